$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying numeric-looking values are stored as text in the sheet
# (numberStoredAsText). Force the data range to remain text so Excel does
# not silently convert the re-assigned values into numbers.
$ws.Range("C2:F8").NumberFormat = "@"

# Row 2
$ws.Range("C2").Value = "0"
$ws.Range("D2").Value = "1"
$ws.Range("E2").Value = "0"
$ws.Range("F2").Value = "0"

# Row 3
$ws.Range("C3").Value = "14"
$ws.Range("D3").Value = "15"
$ws.Range("E3").Value = "1"
$ws.Range("F3").Value = "0"

# Row 4
$ws.Range("C4").Value = "17"
$ws.Range("D4").Value = "14"
$ws.Range("E4").Value = "1"
$ws.Range("F4").Value = "0"

# Row 5
$ws.Range("C5").Value = "66"
$ws.Range("D5").Value = "56"
$ws.Range("E5").Value = "8"
$ws.Range("F5").Value = "2"

# Row 8
$ws.Range("C8").Value = "6"
$ws.Range("D8").Value = "6"
$ws.Range("E8").Value = "0"
$ws.Range("F8").Value = "0"
